$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three new rows before the current "Total" row (row 19), ---
# --- shifting rows 19-35 down to 22-38.                             ---
$ws.Rows.Item(19).Resize(3).Insert()

# Row 19: new section header "ADDITIONS: made on 7/12/2017"
$ws.Cells.Item(19, 1).Value = "ADDITIONS: made on 7/12/2017"
# (style after insert defaults to the row-above's non-header style on col A;
#  fix it to match the other section-header rows, e.g. row 17)
$ws.Cells.Item(17, 1).Copy()
$ws.Cells.Item(19, 1).PasteSpecial(-4122)

# Row 20 / Row 21: new XT30 / XT60 battery connector line items.
$ws.Cells.Item(20, 1).Value = "XT30 Battery Connectors"
$ws.Cells.Item(21, 1).Value = "XT60 Battery Connectors"

$ws.Cells.Item(20, 6).Value = "https://www.amazon.com/Female-Bullet-Connectors-Shrink-Battery/dp/B06ZZSKSJ2/ref=sr_1_3?ie=UTF8&qid=1499883947&sr=8-3&keywords=XT30+battery+connector"
$ws.Cells.Item(21, 2).Value = "10 Pair XT60 Male Female Battery Connectors"
$ws.Cells.Item(20, 2).Value = "10 Pair XT30 Male Female Battery Connectors"
$ws.Cells.Item(21, 6).Value = "https://www.amazon.com/LONMAX-Pairs-Connectors-Battery-Female/dp/B07251HPTQ/ref=sr_1_5?ie=UTF8&qid=1499884651&sr=8-5&keywords=xt60+connectors"

$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = 9.88

$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 6.99

# Assign both subtotal formulas in a single range-level call so the engine
# groups E20:E21 together as one shared formula (matching the existing
# shared-formula pattern used by E12:E18 above them).
$ws.Range("E20:E21").Formula = "=C20*D20"

# Row 22 (previously row 19, shifted by the insert): Total formula needs to
# grow to include the two new line-item rows.
$ws.Cells.Item(22, 5).Formula = "=SUM(E2:E21)"

# --- Update selection to match the saved workbook state ---
$ws.Range("E19").Select()
